# Add a new "performance" column (J) computing FIR / (FFT final with FFT_length = 4 * filter size)
# for every data row, as a single shared formula (J4:J14 share J4's formula; J3 holds its own copy).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J3").Formula = "=B3/H3"
$ws.Range("J4:J14").Formula = "=B4/H4"

# Restore the recorded selection state after the edit.
$ws.Range("F23").Select()
